# Adding the changes we made on may 9th
# Insert two new rows of data at the top (new rows 2-3), shift existing
# rows 2-21 down to rows 4-23, and append eight new rows (24-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.1050096067542932
$ws.Cells.Item(2, 2).Value = -1.756468223065746
$ws.Cells.Item(2, 3).Value = 0.4945203567645989

$ws.Cells.Item(3, 1).Value = -0.0286234012063665
$ws.Cells.Item(3, 2).Value = -0.7998002785809195
$ws.Cells.Item(3, 3).Value = 0.0811297598541999

$ws.Cells.Item(4, 1).Value = 0.03028146071093407
$ws.Cells.Item(4, 2).Value = 0.6503314929349078
$ws.Cells.Item(4, 3).Value = -0.2584614434412545

$ws.Cells.Item(5, 1).Value = 0.6556547600395809
$ws.Cells.Item(5, 2).Value = 3.617322595752007
$ws.Cells.Item(5, 3).Value = -1.325163067603595

$ws.Cells.Item(6, 1).Value = 0.9173419645854408
$ws.Cells.Item(6, 2).Value = 4.483878631981052
$ws.Cells.Item(6, 3).Value = -1.026156997194096

$ws.Cells.Item(7, 1).Value = 0.6556111379545562
$ws.Cells.Item(7, 2).Value = 2.95567157317181
$ws.Cells.Item(7, 3).Value = 0.3850134197546516

$ws.Cells.Item(8, 1).Value = 0.01345462458474067
$ws.Cells.Item(8, 2).Value = 2.573593986277679
$ws.Cells.Item(8, 3).Value = -0.1254206044333366

$ws.Cells.Item(9, 1).Value = -0.0735905529284972
$ws.Cells.Item(9, 2).Value = 1.541839412736652
$ws.Cells.Item(9, 3).Value = -0.5845232478209913

$ws.Cells.Item(10, 1).Value = 0.5347626531610679
$ws.Cells.Item(10, 2).Value = -1.351723280640276
$ws.Cells.Item(10, 3).Value = 0.006906517306153076

$ws.Cells.Item(11, 1).Value = 0.2076723618166787
$ws.Cells.Item(11, 2).Value = -3.982405117579868
$ws.Cells.Item(11, 3).Value = -0.1024508295314652

$ws.Cells.Item(12, 1).Value = -0.6092913156869469
$ws.Cells.Item(12, 2).Value = -4.47834036301593
$ws.Cells.Item(12, 3).Value = 0.3969252789203029

$ws.Cells.Item(13, 1).Value = -0.4407143246154432
$ws.Cells.Item(13, 2).Value = -3.636212723595754
$ws.Cells.Item(13, 3).Value = 0.6680403655889118

$ws.Cells.Item(14, 1).Value = 0.2482107579708095
$ws.Cells.Item(14, 2).Value = -2.94579162889598
$ws.Cells.Item(14, 3).Value = 0.002957710806204883

$ws.Cells.Item(15, 1).Value = 0.1622283559064473
$ws.Cells.Item(15, 2).Value = 0.07419828006199625
$ws.Cells.Item(15, 3).Value = -0.1648120685499541

$ws.Cells.Item(16, 1).Value = 0.2419649730531533
$ws.Cells.Item(16, 2).Value = 4.375038618944129
$ws.Cells.Item(16, 3).Value = -0.3818312956362353

$ws.Cells.Item(17, 1).Value = 0.5893073264433398
$ws.Cells.Item(17, 2).Value = 6.328164937544838
$ws.Cells.Item(17, 3).Value = 0.07646724885823974

$ws.Cells.Item(18, 1).Value = -0.175296502453943
$ws.Cells.Item(18, 2).Value = 2.29715876919883
$ws.Cells.Item(18, 3).Value = 1.448470601013731

$ws.Cells.Item(19, 1).Value = 0.8341801327710269
$ws.Cells.Item(19, 2).Value = 0.753240480714915
$ws.Cells.Item(19, 3).Value = -0.1038439748238527

$ws.Cells.Item(20, 1).Value = -0.1397884144466732
$ws.Cells.Item(20, 2).Value = 0.2036113617371567
$ws.Cells.Item(20, 3).Value = -0.8063265510967792

$ws.Cells.Item(21, 1).Value = -0.7627712436476506
$ws.Cells.Item(21, 2).Value = -1.678813515877241
$ws.Cells.Item(21, 3).Value = -1.093177660387389

$ws.Cells.Item(22, 1).Value = -0.0004051567948546442
$ws.Cells.Item(22, 2).Value = -2.697911312993704
$ws.Cells.Item(22, 3).Value = -0.4782794093599184

$ws.Cells.Item(23, 1).Value = -0.8633800915309378
$ws.Cells.Item(23, 2).Value = -0.3154059344408438
$ws.Cells.Item(23, 3).Value = 0.4824308418497782

$ws.Cells.Item(24, 1).Value = -0.4081483519807154
$ws.Cells.Item(24, 2).Value = -0.6726997543354425
$ws.Cells.Item(24, 3).Value = -0.2190668820118418

$ws.Cells.Item(25, 1).Value = 0.2211332225373814
$ws.Cells.Item(25, 2).Value = 0.241335413285664
$ws.Cells.Item(25, 3).Value = 0.08368853798934378

$ws.Cells.Item(26, 1).Value = 0.06768137718341787
$ws.Cells.Item(26, 2).Value = 0.3379019900244107
$ws.Cells.Item(26, 3).Value = 0.1505034766635118

$ws.Cells.Item(27, 1).Value = 0.07254024853511698
$ws.Cells.Item(27, 2).Value = 0.5556785336562575
$ws.Cells.Item(27, 3).Value = -0.05807583201296457

$ws.Cells.Item(28, 1).Value = 0.1816357883567719
$ws.Cells.Item(28, 2).Value = 0.1322741392923868
$ws.Cells.Item(28, 3).Value = -0.08515337003128903

$ws.Cells.Item(29, 1).Value = -0.02734556931013958
$ws.Cells.Item(29, 2).Value = -0.1169588795425942
$ws.Cells.Item(29, 3).Value = 0.04497027853313797

$ws.Cells.Item(30, 1).Value = -0.02540700723017953
$ws.Cells.Item(30, 2).Value = -0.06986615411481072
$ws.Cells.Item(30, 3).Value = -0.074921377335808

$ws.Cells.Item(31, 1).Value = 0.02237761537639455
$ws.Cells.Item(31, 2).Value = -0.07008743807863513
$ws.Cells.Item(31, 3).Value = -0.003453258577050004

